# "Generate Report for Handoff"
# The localization run moved from "In Translation" to "Ready for handoff":
#   - Overview sheet: both locale status columns (E/F) and the HO xliff
#     generate timestamp (G) are refreshed.
#   - zh-cn sheet: Status column (C) and Latest Handoff Datetime (H).
#   - de-de sheet: Status column (C) and Latest Handoff Datetime (H).
# The new status text is longer than the old one, so the Status columns
# widen to fit the new content (matches the wider columns seen after a
# handoff report regeneration).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newWidth = 16.3333333333333   # widened Status column to fit "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-25 00:36:39"
$wsOverview.Range("E1").ColumnWidth = $newWidth
$wsOverview.Range("F1").ColumnWidth = $newWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-25 00:36:35"
$wsZhCn.Range("C1").ColumnWidth = $newWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-25 00:36:39"
$wsDeDe.Range("C1").ColumnWidth = $newWidth
